$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new quarterly columns before column D. Everything that used
#    to live in D:K (8 quarters) shifts right to F:M; Excel carries the
#    values/formats/shared-strings along automatically.
# ---------------------------------------------------------------------------
$ws.Range("D1:E1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 2. The freshly inserted D:E columns pick up column C's format by default.
#    Copy the number format from column F (the old column D, now shifted)
#    into D and E so date rows keep the date style and value rows keep the
#    numeric style, matching what Excel does when you insert+fill columns.
# ---------------------------------------------------------------------------
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Label-only rows (5, 6, 37, 79) should not gain D/E cells at all - undo the
# formatting paste there.
$ws.Range("D5:E6").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# ---------------------------------------------------------------------------
# 3. Fill the two new quarter columns (D = 2018-12-31, E = 2018-09-30) with
#    the newly reported data for every row that carries data.
#    Format per line: "<row>|<D value>|<E value>"; "NA" means the literal
#    text "NA"; "$null" means leave the cell as-is (no value in this row).
# ---------------------------------------------------------------------------
$newData = @'
7|43465|43373
8|3712000|3504000
9|NA|NA
10|NA|NA
11|$null|$null
12|NA|NA
13|0|0
14|0|0
15|0|0
16|$null|$null
17|3091000|2963000
18|621000|541000
19|$null|$null
20|-305000|-87000
21|439000|578000
22|92000|69000
23|224000|385000
24|59000|120000
25|0|0
26|165000|265000
27|159000|262000
28|0|0
29|-6000|14000
30|0|0
31|0|0
32|305000|87000
33|153000|276000
34|0|0
35|153000|276000
38|43465|43373
39|$null|$null
40|$null|$null
41|1066000|951000
42|0|0
43|4317000|4476000
44|0|0
45|551000|539000
46|5934000|5966000
47|0|0
48|701000|707000
49|11036000|10764000
50|0|0
51|0|0
52|3907000|3692000
53|0|0
54|21578000|21129000
55|$null|$null
56|$null|$null
57|2234000|2293000
58|314000|638000
59|2376000|1796000
60|4924000|4727000
61|5510000|5512000
62|3560000|3352000
63|0|0
64|0|0
65|0|0
66|14067000|13676000
67|$null|$null
68|0|0
69|0|0
70|0|0
71|0|0
72|14347000|14196000
73|0|0
74|0|0
75|0|0
76|7511000|7453000
77|0|0
80|43465|43373
81|153000|276000
82|$null|$null
83|123000|124000
84|0|0
85|0|0
86|0|0
87|0|0
88|0|0
89|1109000|906000
90|$null|$null
91|-92000|-87000
92|0|0
93|0|0
94|-337000|-473000
95|$null|$null
96|-213000|-211000
97|0|0
98|0|0
99|0|0
100|-563000|-515000
101|-94000|-3000
102|115000|-85000
'@

foreach ($line in ($newData -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $r = [int]$parts[0]
    $dtok = $parts[1]
    $etok = $parts[2]

    if ($dtok -ne '$null') {
        if ($dtok -eq 'NA') {
            $ws.Cells.Item($r, 4).Value2 = "NA"
        } else {
            $ws.Cells.Item($r, 4).Value2 = [double]$dtok
        }
    }
    if ($etok -ne '$null') {
        if ($etok -eq 'NA') {
            $ws.Cells.Item($r, 5).Value2 = "NA"
        } else {
            $ws.Cells.Item($r, 5).Value2 = [double]$etok
        }
    }
}

# ---------------------------------------------------------------------------
# 4. A handful of historical quarters were also restated/corrected alongside
#    the refresh. These land in the shifted columns (old D->F, old F->H,
#    old G->I) after the column insert above.
# ---------------------------------------------------------------------------
$corrections = @'
17|H|3015000
17|I|2806000
18|H|670000
18|I|535000
20|H|31000
20|I|62000
24|F|183000
26|F|536000
27|F|531000
29|F|0
32|H|-31000
32|I|-62000
'@

$colIndex = @{ 'F' = 6; 'G' = 7; 'H' = 8; 'I' = 9; 'J' = 10; 'K' = 11; 'L' = 12; 'M' = 13 }

foreach ($line in ($corrections -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $r = [int]$parts[0]
    $col = $colIndex[$parts[1]]
    $val = [double]$parts[2]
    $ws.Cells.Item($r, $col).Value2 = $val
}

Write-Output "done"
